# Insert two new rows at 426-427 (existing rows 426-453 shift down to 428-455)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("426:427").Insert()

# Row 426: new "Primera" record dated 44714
$ws.Cells.Item(426,1).Value  = 4
$ws.Cells.Item(426,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(426,3).Value  = "Los Lagos"
$ws.Cells.Item(426,4).Value  = 44714
$ws.Cells.Item(426,5).Value  = 10
$ws.Cells.Item(426,6).Value  = "Fruta"
$ws.Cells.Item(426,7).Value  = 100106
$ws.Cells.Item(426,8).Value  = "Oleaginosos"
$ws.Cells.Item(426,9).Value  = 100106002
$ws.Cells.Item(426,10).Value = "Palta"
$ws.Cells.Item(426,11).Value = "Hass"
$ws.Cells.Item(426,12).Value = "Primera"
$ws.Cells.Item(426,13).Value = 200
$ws.Cells.Item(426,14).Value = 4300
$ws.Cells.Item(426,15).Value = 4400
$ws.Cells.Item(426,16).Value = 4350
$ws.Cells.Item(426,17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(426,18).Value = "Provincia de Quillota"
$ws.Cells.Item(426,19).Value = 4350
$ws.Cells.Item(426,20).Value = 1

# Row 427: new "Segunda" record dated 44714
$ws.Cells.Item(427,1).Value  = 4
$ws.Cells.Item(427,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427,3).Value  = "Los Lagos"
$ws.Cells.Item(427,4).Value  = 44714
$ws.Cells.Item(427,5).Value  = 10
$ws.Cells.Item(427,6).Value  = "Fruta"
$ws.Cells.Item(427,7).Value  = 100106
$ws.Cells.Item(427,8).Value  = "Oleaginosos"
$ws.Cells.Item(427,9).Value  = 100106002
$ws.Cells.Item(427,10).Value = "Palta"
$ws.Cells.Item(427,11).Value = "Hass"
$ws.Cells.Item(427,12).Value = "Segunda"
$ws.Cells.Item(427,13).Value = 100
$ws.Cells.Item(427,14).Value = 4000
$ws.Cells.Item(427,15).Value = 4000
$ws.Cells.Item(427,16).Value = 4000
$ws.Cells.Item(427,17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(427,18).Value = "Provincia de Quillota"
$ws.Cells.Item(427,19).Value = 4000
$ws.Cells.Item(427,20).Value = 1
